# Update on 20250810 part 14
# Shorten the script list in cell C4 of sheet "地方台JS脚本" (sheet1):
# remove the second line "Shanghai_setv2.php", keeping only "Shanghai_setv.php".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("地方台JS脚本")

$ws.Range("C4").Value = "Shanghai_setv.php"

# The row no longer needs the extra height required to show two wrapped
# lines of text, so let Excel recompute the row height automatically.
$ws.Rows.Item(4).AutoFit()
